$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert a brand-new "2022-Q1" sheet right after "2021-Q4" and before the
#    "总计" (totals) sheet, following the same layout as the other per-
#    quarter sheets (fund code / name / size / position / rank table).
# ---------------------------------------------------------------------------
$q4 = $wb.Worksheets.Item("2021-Q4")
$newQ = $wb.Worksheets.Add($null, $q4)
$newQ.Name = "2022-Q1"

# Pull over the header-row and index-column formatting (bold + border +
# centered, style used by every other quarter sheet) from "2021-Q4" so the
# new sheet's cellXfs line up with the rest of the workbook.
$q4.Range("B1:H1").Copy()
$newQ.Range("B1:H1").PasteSpecial(-4122)
$q4.Range("A2:A4").Copy()
$newQ.Range("A2:A4").PasteSpecial(-4122)

$newQ.Range("B1").Value = "基金代码"
$newQ.Range("C1").Value = "基金名称"
$newQ.Range("D1").Value = "基金规模"
$newQ.Range("E1").Value = "股票总仓位"
$newQ.Range("F1").Value = "仓位占比"
$newQ.Range("G1").Value = "持有市值(亿元)"
$newQ.Range("H1").Value = "仓位排名"

$newQ.Range("A2").Value = 0
$newQ.Range("A3").Value = 1
$newQ.Range("A4").Value = 2

# B:G on the data rows are stored as text in the source data, so force text
# formatting before assigning, then strip the format again so no stray
# numFmt-only style lingers on the cells.
$dataText = $newQ.Range("B2:G4")
$dataText.NumberFormat = "@"

$newQ.Range("B2").Value = "257010"
$newQ.Range("C2").Value = "国联安小盘精选混合"
$newQ.Range("D2").Value = "9.15"
$newQ.Range("E2").Value = "74.36"
$newQ.Range("F2").Value = "6.05"
$newQ.Range("G2").Value = "0.5536"
$newQ.Range("H2").Value = 1

$newQ.Range("B3").Value = "002367"
$newQ.Range("C3").Value = "国联安安稳灵活配置混合"
$newQ.Range("D3").Value = "2.32"
$newQ.Range("E3").Value = "33.99"
$newQ.Range("F3").Value = "3.89"
$newQ.Range("G3").Value = "0.0902"
$newQ.Range("H3").Value = 2

$newQ.Range("B4").Value = "006138"
$newQ.Range("C4").Value = "国联安价值优选股票"
$newQ.Range("D4").Value = "0.60"
$newQ.Range("E4").Value = "93.30"
$newQ.Range("F4").Value = "5.33"
$newQ.Range("G4").Value = "0.0320"
$newQ.Range("H4").Value = 3

$dataText.ClearFormats()

# ---------------------------------------------------------------------------
# 2. Update the "总计" (totals) sheet: add a new top data row for 2022-Q1,
#    pushing the existing quarters down and renumbering the running index
#    in column A.
# ---------------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item("总计")

# Push the 5 existing data rows down one slot (bottom-up so nothing is
# clobbered before it's moved). Using Cut/paste-by-move keeps each row's
# original style/type intact and avoids minting any new cellXfs entries
# (unlike Rows.Insert, which drags in formatting from the row above).
$wsTotal.Range("A6:D6").Cut($wsTotal.Range("A7:D7"))
$wsTotal.Range("A5:D5").Cut($wsTotal.Range("A6:D6"))
$wsTotal.Range("A4:D4").Cut($wsTotal.Range("A5:D5"))
$wsTotal.Range("A3:D3").Cut($wsTotal.Range("A4:D4"))
$wsTotal.Range("A2:D2").Cut($wsTotal.Range("A3:D3"))

# Row 2 (now vacated) becomes the new 2022-Q1 entry. A2 already carries the
# index column's bold/border/center style since it wasn't touched by the
# cuts above.
$wsTotal.Range("A2").Value = 0
$wsTotal.Range("B2").NumberFormat = "@"
$wsTotal.Range("B2").Value = "2022-Q1"
$wsTotal.Range("B2").ClearFormats()
$wsTotal.Range("C2").Value = 3
$wsTotal.Range("D2").Value = 0.68

# Renumber the running index in column A for the rows that shifted down.
$wsTotal.Range("A3").Value = 1
$wsTotal.Range("A4").Value = 2
$wsTotal.Range("A5").Value = 3
$wsTotal.Range("A6").Value = 4
$wsTotal.Range("A7").Value = 5
